$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed cells retain their original text (inline string) type
# by pre-formatting them as Text before assigning the new values.
$changedCells = @("D2","G2","D3","G3","D4","G4","D5","G5","D6","G6","G7","D8","G8","D9","G9","B10","C10","D10","E10","G10","B11","C11","D11","E11","G11","B12","C12","D12","E12","G12","B13","C13","D13","E13","G13","B14","C14","D14","E14","G14","B15","C15","D15","E15","G15","B16","C16","D16","E16","G16","B17","C17","D17","E17","G17","B18","C18","D18","E18","G18","D19","G19","G20","D21","G21","G22","D23","G23","D24","G24","D25","G25","G26","E27","G27","G28","G29","G30","G31","G32","G33","G34","G35","G36","G37","G38","G39","D40","G40","G41","B42","C42","D42","E42","G42","B43","C43","D43","E43","G43","D44","G44","D45","G45","G46","E47","G47","D48","G48","G49","G50","G51")
foreach ($addr in $changedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from the latest coinranking.com snapshot
$ws.Range('D2').Value = '243.41'
$ws.Range('G2').Value = '4'
$ws.Range('D3').Value = '23.15'
$ws.Range('G3').Value = '4'
$ws.Range('D4').Value = '5.405'
$ws.Range('G4').Value = '4'
$ws.Range('D5').Value = '0.05980'
$ws.Range('G5').Value = '4'
$ws.Range('D6').Value = '3.429'
$ws.Range('G6').Value = '4'
$ws.Range('G7').Value = '4'
$ws.Range('D8').Value = '0.8088'
$ws.Range('G8').Value = '4'
$ws.Range('D9').Value = '0.9190'
$ws.Range('G9').Value = '4'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').Value = '0.01114'
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('G10').Value = '4'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '0.1428'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('G11').Value = '4'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '0.07437'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('G12').Value = '4'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D13').Value = '0.03304'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('G13').Value = '4'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').Value = '0.03066'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('G14').Value = '4'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = '0.09355'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('G15').Value = '4'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').Value = '3.853'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('G16').Value = '4'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').Value = '0.001585'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('G17').Value = '4'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').Value = '0.04700'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('G18').Value = '4'
$ws.Range('D19').Value = '0.005866'
$ws.Range('G19').Value = '4'
$ws.Range('G20').Value = '4'
$ws.Range('D21').Value = '0.004884'
$ws.Range('G21').Value = '4'
$ws.Range('G22').Value = '4'
$ws.Range('D23').Value = '3.570'
$ws.Range('G23').Value = '4'
$ws.Range('D24').Value = '2.138'
$ws.Range('G24').Value = '4'
$ws.Range('D25').Value = '0.3233'
$ws.Range('G25').Value = '4'
$ws.Range('G26').Value = '4'
$ws.Range('E27').Value = '26UpBotsUBXTWorstin24h'
$ws.Range('G27').Value = '4'
$ws.Range('G28').Value = '4'
$ws.Range('G29').Value = '4'
$ws.Range('G30').Value = '4'
$ws.Range('G31').Value = '4'
$ws.Range('G32').Value = '4'
$ws.Range('G33').Value = '4'
$ws.Range('G34').Value = '4'
$ws.Range('G35').Value = '4'
$ws.Range('G36').Value = '4'
$ws.Range('G37').Value = '4'
$ws.Range('G38').Value = '4'
$ws.Range('G39').Value = '4'
$ws.Range('D40').Value = '0.03970'
$ws.Range('G40').Value = '4'
$ws.Range('G41').Value = '4'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Value = '0.1077'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('G42').Value = '4'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').Value = '0.003100'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('G43').Value = '4'
$ws.Range('D44').Value = '0.009836'
$ws.Range('G44').Value = '4'
$ws.Range('D45').Value = '0.00005070'
$ws.Range('G45').Value = '4'
$ws.Range('G46').Value = '4'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOIN'
$ws.Range('G47').Value = '4'
$ws.Range('D48').Value = '0.002427'
$ws.Range('G48').Value = '4'
$ws.Range('G49').Value = '4'
$ws.Range('G50').Value = '4'
$ws.Range('G51').Value = '4'
